$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 5, 7, 9, 13, 14, 16, 20, 22, 26)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "done"
}

$ws.Range("B19").Select()
